$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold font, border, centered/top alignment) from H1 into I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF data values for rows 2-84
$iVals = @(8,9,9,9,9,9,8,9,8,9,9,9,8,9,9,9,9,8,7,6,6,9,7,8,6,7,7,7,9,6,7,6,6,7,7,8,6,7,8,9,7,7,8,6,7,7,8,8,8,8,9,7,6,7,7,6,7,8,8,6,8,7,1,7,5,9,6,8,8,7,7,8,7,6,9,7,7,6,6,5,7,4,4)
$jVals = @(8,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,8,7,7,6,9,7,8,6,7,7,7,9,7,7,7,7,7,8,8,7,7,8,9,7,7,8,7,7,7,8,8,8,8,9,7,6,7,8,6,7,8,8,6,8,7,1,7,5,9,6,8,8,7,7,8,8,6,9,7,7,6,6,5,7,4,4)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}

Write-Output "done"
